$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.085.21"
$ws.Range("E2").Value = "  +4.31%  "
$ws.Range("D3").Value = "2.664.88"
$ws.Range("E3").Value = "  +7.32%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +8.44%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "326.81"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +2.09%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +3.75%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "41.39"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +6.39%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "20.18"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("E14").Value = "  +4.42%  "
$ws.Range("D15").Value = "3.078.41"
$ws.Range("E15").Value = "  +7.11%  "
$ws.Range("D16").Value = "2.663.96"
$ws.Range("E16").Value = "  +7.04%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.881"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +6.74%  "
$ws.Range("D18").Value = "50.025.06"
$ws.Range("E18").Value = "  +4.39%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.35"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  +3.35%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "72.65"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "278.29"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +2.60%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.59"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.50%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "26.96"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +5.13%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "10.06"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +3.97%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "36.79"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +6.57%  "
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("E33").Value = "  +5.01%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "19.75"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("E35").Value = "  +6.56%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.08"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +11.84%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +8.20%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.14"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +9.77%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "125.40"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "22.18"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("E44").Value = "  +5.48%  "
$ws.Range("D45").Value = "2.112.23"
$ws.Range("E45").Value = "  +5.58%  "
$ws.Range("E46").Value = "  +6.38%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +14.23%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +4.88%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.11"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.35%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "5.38"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +4.36%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "59.66"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +5.91%  "
